$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = 2024
$ws.Range("B2").Value = "DEC"
$ws.Range("C2").Value = "31/12-01/12"
$ws.Range("D2").Value = "-"
$ws.Range("E2").Value = "Akurana"

$ws.Range("F2").Value = 24.11539757412399
$ws.Range("G2").Value = 25.38117699910153
$ws.Range("H2").Value = 21.51930592991914
$ws.Range("I2").Value = 0.6678571428571429
$ws.Range("J2").Value = 1.766531895777179
$ws.Range("K2").Value = 8.521585804132972
$ws.Range("L2").Value = 186.3287286612758
$ws.Range("M2").Value = 5.361668912848159
$ws.Range("N2").Value = 1.301403863432165
$ws.Range("O2").Value = 9.480491913746631
$ws.Range("P2").Value = 55.53574797843666
$ws.Range("Q2").Value = 2459.035747978436
$ws.Range("R2").Value = 960.4171759155245
$ws.Range("S2").Value = 86.86769991015274
$ws.Range("T2").Value = 1.055143755615454
$ws.Range("U2").Value = 131.0022012578617
$ws.Range("V2").Value = 24.12546618737362
$ws.Range("W2").Value = 85.02662323073467
$ws.Range("X2").Value = 1
$ws.Range("Y2").Value = 24.86379465288699
$ws.Range("Z2").Value = 21.41085149404628
$ws.Range("AA2").Value = 960.4171759155245
$ws.Range("AB2").Value = 1
$ws.Range("AC2").Value = 186.7412398921833
$ws.Range("AD2").Value = 0.5443059299191375

$ws.Range("AE2:AJ2").ClearContents()
